# "Put Names and Courses on the home page"
#
# Classes sheet:
#   - E1 header "Credit hours needed" -> "GA hours needed"
#   - Remove the "Class nomination" / "Professor Nominations" columns (K, L)
#     (their content moves to the Ga's sheet, see below)
#   - Add a new row for "Full Department"
#
# Ga's sheet:
#   - F1 header "GA nomination" -> "Hours", and the Hours values shift from
#     column G into column F
#   - New "Class nomination" / "Professor Nominations" columns (G, H) added
#     at the end, taking over from the Classes sheet

$wb = $excel.ActiveWorkbook

$classes = $wb.Worksheets.Item("Classes")
$gas = $wb.Worksheets.Item("Ga's")

# --- Classes sheet -------------------------------------------------------

# Rename "Credit hours needed" to "GA hours needed"
$classes.Range("E1").Value = "GA hours needed"

# Drop the old "Class nomination" / "Professor Nominations" columns - they
# move over to the Ga's sheet
$classes.Range("K1:L1").ClearContents()

# New row for the whole department
$classes.Range("B7").Value = "Full Department"

$classes.Range("F18").Select()

# --- Ga's sheet ------------------------------------------------------------

# "GA nomination" column becomes "Hours", sliding the existing Hours values
# from column G back into column F
$gas.Range("F1").Value = "Hours"
$gas.Range("F2").Value = $gas.Range("G2").Value2
$gas.Range("F3").Value = $gas.Range("G3").Value2
$gas.Range("F4").Value = $gas.Range("G4").Value2
$gas.Range("F5").Value = $gas.Range("G5").Value2
$gas.Range("F6").Value = $gas.Range("G6").Value2
$gas.Range("F7").Value = $gas.Range("G7").Value2
$gas.Range("G2:G7").ClearContents()

# New columns for class/professor nominations
$gas.Range("G1").Value = "Class nomination"
$gas.Range("H1").Value = "Professor Nominations "
$gas.Range("H2").Value = "Sen"

$gas.Range("D1").Select()

$classes.Activate()
